# chore: adapt column header formatting to respective input file names
#
# The sheet holds a side-by-side AHB diff: the left block of columns
# (A:J) described the "old" format version and the right block (L:U)
# the "new" one. Rename those header suffixes to the concrete format
# versions being compared (FV2210 / FV2304), wrap the data range in a
# native Excel Table so the headers carry filter buttons, and freeze
# the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1): "_old" -> "_FV2210", "_new" -> "_FV2304" ---
$headerCols = @("A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1","Q1","R1","S1","T1","U1")
$headerNames = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i]).Value = $headerNames[$i]
}

# --- 2. Turn the header + data range into a native table (adds the autofilter) ---
$dataRange = $ws.Range("A1:U94")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze panes below the header row (row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
